# "call not support '$' replace '.'"
# Slides 2-4 each contain:
#   - two ellipses inside a nested group reading "SystemA $FuncN"
#   - a textbox reading "Loading Name $ Api " / "...’$ ‘..."
# Replace the literal '$' call-symbol with '.' everywhere it is used to
# denote the call syntax, without disturbing any other run formatting.

$p = $ppt.ActivePresentation

# Unicode right/left single quotation marks used around the '$' token in
# the explanatory sentence ( ’$ ‘ -> ’. ‘ ).
$rsquo = [char]0x2019
$lsquo = [char]0x2018

for ($si = 2; $si -le 4; $si++) {
    $s = $p.Slides.Item($si)

    # "그룹 7" -> (flattened) "타원 10" / "타원 11"
    $grp = $s.Shapes.Item(4)
    $ellipse1 = $grp.GroupItems.Item(2)   # SystemA $Func1
    $ellipse2 = $grp.GroupItems.Item(3)   # SystemA $Func2

    $ellipse1.TextFrame.TextRange.Text = "SystemA .Func1"
    $ellipse2.TextFrame.TextRange.Text = "SystemA .Func2"

    # "TextBox 9"
    $tb = $s.Shapes.Item(5)
    $tr = $tb.TextFrame.TextRange

    # Paragraph 1, single run: "Loading Name $ Api " -> "Loading Name . Api "
    $tr.Characters(1, 20).Text = "Loading Name . Api "

    # Paragraph 2, 4th run: "<U+2019>$ <U+2018> " -> "<U+2019>. <U+2018> "
    $tr.Characters(31, 5).Text = $rsquo + ". " + $lsquo + " "
}
